$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# Row 3
$ws.Range("E3").Value = 10.78
$ws.Range("F3").Value = 10.34

# Row 4
$ws.Range("E4").Value = 10.6

# Row 5
$ws.Range("C5").Value = 9.18
$ws.Range("D5").Value = 9.4
$ws.Range("F5").Value = 10.17
$ws.Range("G5").Value = 9.65

# Row 6
$ws.Range("C6").Value = 9.66
$ws.Range("E6").Value = 9.83
$ws.Range("G6").Value = 10.42
$ws.Range("H6").Value = 10.32

# Row 7
$ws.Range("E7").Value = 10.35
$ws.Range("F7").Value = 9.58
$ws.Range("J7").Value = 9.88

# Row 8
$ws.Range("F8").Value = 9.68

# Row 10
$ws.Range("G10").Value = 10.12
